$wb = $excel.ActiveWorkbook

# ALC row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 166.78572
$ws.Range("I12").Value = 166.78572
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 166.78572
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = 3.214280000000002

# ALC row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 246.33333
$ws.Range("I18").Value = 246.33333
$ws.Range("K18").Value = 246.33333
$ws.Range("M18").Value = 37.66667000000001

# ALC row 47
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 17500
$ws.Range("J47").Value = 17500
$ws.Range("L47").Value = 17500
$ws.Range("N47").Value = -19444

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6749
$ws.Range("I62").Value = 6823.9
$ws.Range("K62").Value = 6823.9
$ws.Range("M62").Value = -6199.9

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3733.3333
$ws.Range("J64").Value = 4000
$ws.Range("L64").Value = 4000
$ws.Range("N64").Value = -4496

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 6749
$ws.Range("I65").Value = 6823.9
$ws.Range("K65").Value = 34119.5
$ws.Range("M65").Value = -30999.5

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3733.3333
$ws.Range("J67").Value = 4000
$ws.Range("L67").Value = 4000
$ws.Range("N67").Value = -5716

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1198
$ws.Range("I86").Value = 946.6667
$ws.Range("J86").Value = 1499.6
$ws.Range("K86").Value = 946.6667
$ws.Range("L86").Value = 1499.6
$ws.Range("M86").Value = 176.3333
$ws.Range("N86").Value = -3745.6

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 733.1539
$ws.Range("J88").Value = 697.375
$ws.Range("L88").Value = 697.375
$ws.Range("N88").Value = -1509.375

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1198
$ws.Range("I89").Value = 946.6667
$ws.Range("J89").Value = 1499.6
$ws.Range("K89").Value = 4733.3335
$ws.Range("L89").Value = 7498
$ws.Range("M89").Value = 882.6665000000003
$ws.Range("N89").Value = -18730

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 733.1539
$ws.Range("J91").Value = 697.375
$ws.Range("L91").Value = 697.375
$ws.Range("N91").Value = -3505.375

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1023.05
$ws.Range("I107").Value = 1073.125
$ws.Range("J107").Value = 822.75
$ws.Range("K107").Value = 1073.125
$ws.Range("L107").Value = 822.75
$ws.Range("M107").Value = 846.875
$ws.Range("N107").Value = -4662.75

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 7173
$ws.Range("I113").Value = 5008.1665
$ws.Range("K113").Value = 5008.1665
$ws.Range("M113").Value = -1754.1665

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2675.923
$ws.Range("I125").Value = 3466
$ws.Range("J125").Value = 2324.7778
$ws.Range("K125").Value = 31194
$ws.Range("L125").Value = 20923.0002
$ws.Range("M125").Value = -28734
$ws.Range("N125").Value = -25843.0002

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2599.75
$ws.Range("I132").Value = 1799.8334
$ws.Range("K132").Value = 5399.5002
$ws.Range("M132").Value = -2869.5002

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4591.2915
$ws.Range("I138").Value = 4710.6
$ws.Range("K138").Value = 14131.8
$ws.Range("M138").Value = -8991.800000000001

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4317.6665
$ws.Range("I141").Value = 4357.5
$ws.Range("J141").Value = 3999
$ws.Range("K141").Value = 13072.5
$ws.Range("L141").Value = 11997
$ws.Range("M141").Value = -7892.5
$ws.Range("N141").Value = -22357

# ARM row 26
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").ClearContents()
$ws.Range("N26").Value = 0

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1924.381
$ws.Range("I32").Value = 1483.5172
$ws.Range("K32").Value = 1483.5172
$ws.Range("M32").Value = -1196.5172

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5532.136
$ws.Range("I61").Value = 3929.4707
$ws.Range("K61").Value = 3929.4707
$ws.Range("M61").Value = -3717.4707

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3585.544
$ws.Range("I74").Value = 2952.0303
$ws.Range("J74").Value = 4456.625
$ws.Range("K74").Value = 2952.0303
$ws.Range("L74").Value = 4456.625
$ws.Range("M74").Value = -2078.0303
$ws.Range("N74").Value = -6204.625

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3585.544
$ws.Range("I77").Value = 2952.0303
$ws.Range("J77").Value = 4456.625
$ws.Range("K77").Value = 14760.1515
$ws.Range("L77").Value = 22283.125
$ws.Range("M77").Value = -10392.1515
$ws.Range("N77").Value = -31019.125

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 6117.7144
$ws.Range("I110").Value = 5982.2
$ws.Range("J110").Value = 6456.5
$ws.Range("K110").Value = 5982.2
$ws.Range("L110").Value = 6456.5
$ws.Range("M110").Value = -3937.2
$ws.Range("N110").Value = -10546.5

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5504.7617
$ws.Range("I132").Value = 5732.5
$ws.Range("K132").Value = 17197.5
$ws.Range("M132").Value = -14667.5

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5532.136
$ws.Range("I136").Value = 3929.4707
$ws.Range("K136").Value = 11788.4121
$ws.Range("M136").Value = -9238.4121

# BSM row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 451.2857
$ws.Range("I80").Value = 272
$ws.Range("J80").Value = 899.5
$ws.Range("K80").Value = 272
$ws.Range("L80").Value = 899.5
$ws.Range("M80").Value = 726
$ws.Range("N80").Value = -2895.5

# BSM row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 25906.043
$ws.Range("J82").Value = 40172.23
$ws.Range("L82").Value = 40172.23
$ws.Range("N82").Value = -40938.23

# BSM row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 451.2857
$ws.Range("I83").Value = 272
$ws.Range("J83").Value = 899.5
$ws.Range("K83").Value = 1360
$ws.Range("L83").Value = 4497.5
$ws.Range("M83").Value = 3632
$ws.Range("N83").Value = -14481.5

# BSM row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 25906.043
$ws.Range("J85").Value = 40172.23
$ws.Range("L85").Value = 40172.23
$ws.Range("N85").Value = -42824.23

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1527.069
$ws.Range("J94").Value = 1281.25
$ws.Range("L94").Value = 1281.25
$ws.Range("N94").Value = -2183.25

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7722.5405
$ws.Range("I134").Value = 4904.3125
$ws.Range("K134").Value = 14712.9375
$ws.Range("M134").Value = -12177.9375

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 906.1667
$ws.Range("I16").Value = 906.1667
$ws.Range("K16").Value = 906.1667
$ws.Range("M16").Value = -619.1667

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2937.4314
$ws.Range("J31").Value = 3403.725
$ws.Range("L31").Value = 3403.725
$ws.Range("N31").Value = -3993.725

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2937.4314
$ws.Range("J34").Value = 3403.725
$ws.Range("L34").Value = 3403.725
$ws.Range("N34").Value = -3807.725

# CRP row 56
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 29998
$ws.Range("J56").Value = 29998
$ws.Range("L56").Value = 29998
$ws.Range("N56").Value = -31688

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5353.909
$ws.Range("I58").Value = 1814.091
$ws.Range("J58").Value = 8893.727999999999
$ws.Range("K58").Value = 1814.091
$ws.Range("L58").Value = 8893.727999999999
$ws.Range("M58").Value = -1611.091
$ws.Range("N58").Value = -9299.727999999999

# CRP row 59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 36163
$ws.Range("J59").Value = 32532.777
$ws.Range("L59").Value = 32532.777
$ws.Range("N59").Value = -34822.777

# CRP row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 34995.5
$ws.Range("I60").Value = 20000
$ws.Range("K60").Value = 20000
$ws.Range("M60").Value = -19489

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3100

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3100

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10643.774
$ws.Range("I99").Value = 5349
$ws.Range("K99").Value = 5349
$ws.Range("M99").Value = -3851

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 906.1667
$ws.Range("I113").Value = 906.1667
$ws.Range("K113").Value = 906.1667
$ws.Range("M113").Value = 1263.8333

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1747
$ws.Range("I122").Value = 1494
$ws.Range("K122").Value = 4482
$ws.Range("M122").Value = -2032

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 10643.774
$ws.Range("I126").Value = 5349
$ws.Range("K126").Value = 16047
$ws.Range("M126").Value = -13577

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 26238.291
$ws.Range("I132").Value = 18121.205
$ws.Range("J132").Value = 39380.24
$ws.Range("K132").Value = 54363.61500000001
$ws.Range("L132").Value = 118140.72
$ws.Range("M132").Value = -51833.61500000001
$ws.Range("N132").Value = -123200.72

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5016.5557
$ws.Range("I134").Value = 4230.2383
$ws.Range("J134").Value = 6117.4
$ws.Range("K134").Value = 12690.7149
$ws.Range("L134").Value = 18352.2
$ws.Range("M134").Value = -10155.7149
$ws.Range("N134").Value = -23422.2

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5353.909
$ws.Range("I136").Value = 1814.091
$ws.Range("J136").Value = 8893.727999999999
$ws.Range("K136").Value = 5442.272999999999
$ws.Range("L136").Value = 26681.184
$ws.Range("M136").Value = -2892.272999999999
$ws.Range("N136").Value = -31781.184

# CRP row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 183240.2
$ws.Range("I141").Value = 83942
$ws.Range("K141").Value = 83942
$ws.Range("M141").Value = -78762

# CUL row 9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1000
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1301.6522
$ws.Range("J131").Value = 1301.7727
$ws.Range("L131").Value = 3905.3181
$ws.Range("N131").Value = -13985.3181

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3387.7334
$ws.Range("J137").Value = 3646.5715
$ws.Range("L137").Value = 10939.7145
$ws.Range("N137").Value = -21139.7145

# GSM row 23
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 4250
$ws.Range("J23").Value = 4250
$ws.Range("L23").Value = 4250
$ws.Range("N23").Value = -4696

# GSM row 55
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").ClearContents()
$ws.Range("N55").Value = 0

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3262.125
$ws.Range("I102").Value = 3540.4546
$ws.Range("J102").Value = 2649.8
$ws.Range("K102").Value = 3540.4546
$ws.Range("L102").Value = 2649.8
$ws.Range("M102").Value = -1918.4546
$ws.Range("N102").Value = -5893.8

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 224266.56
$ws.Range("I113").Value = 287628.44
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 287628.44
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = -285458.44
$ws.Range("N113").Value = -6840

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3583.7856
$ws.Range("I126").Value = 3368.25
$ws.Range("J126").Value = 3871.1667
$ws.Range("K126").Value = 10104.75
$ws.Range("L126").Value = 11613.5001
$ws.Range("M126").Value = -7634.75
$ws.Range("N126").Value = -16553.5001

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 15843.5
$ws.Range("I132").Value = 18299.666
$ws.Range("K132").Value = 54898.99800000001
$ws.Range("M132").Value = -52368.99800000001

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11173.588
$ws.Range("I7").Value = 14953.909
$ws.Range("J7").Value = 4243
$ws.Range("K7").Value = 14953.909
$ws.Range("L7").Value = 4243
$ws.Range("M7").Value = -14841.909
$ws.Range("N7").Value = -4467

# LTW row 20
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 89103.75
$ws.Range("I20").Value = 1000
$ws.Range("J20").Value = 92934.35000000001
$ws.Range("K20").Value = 1000
$ws.Range("L20").Value = 92934.35000000001
$ws.Range("M20").Value = -774
$ws.Range("N20").Value = -93386.35000000001

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3749.25
$ws.Range("I40").Value = 3332.6667
$ws.Range("K40").Value = 3332.6667
$ws.Range("M40").Value = -3196.6667

# LTW row 45
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 26693.666
$ws.Range("I45").Value = 20041
$ws.Range("K45").Value = 20041
$ws.Range("M45").Value = -19634

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1259.6364
$ws.Range("I46").Value = 790.2857
$ws.Range("J46").Value = 1605.4736
$ws.Range("K46").Value = 790.2857
$ws.Range("L46").Value = 1605.4736
$ws.Range("M46").Value = -602.2857
$ws.Range("N46").Value = -1981.4736

# LTW row 47
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 25999
$ws.Range("I47").Value = 25999
$ws.Range("K47").Value = 25999
$ws.Range("M47").Value = -25509

# LTW row 48
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 39999
$ws.Range("J48").Value = 39999
$ws.Range("L48").Value = 39999
$ws.Range("N48").Value = -41321

# LTW row 52
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H52").Value = 25999
$ws.Range("I52").Value = 25999
$ws.Range("K52").Value = 25999
$ws.Range("M52").Value = -25766

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 11173.588
$ws.Range("I126").Value = 14953.909
$ws.Range("J126").Value = 4243
$ws.Range("K126").Value = 44861.727
$ws.Range("L126").Value = 12729
$ws.Range("M126").Value = -42391.727
$ws.Range("N126").Value = -17669

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3630.681
$ws.Range("I136").Value = 3535.6956
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 10607.0868
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -8057.086800000001
$ws.Range("N136").Value = -29100

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5196.5415
$ws.Range("I122").Value = 2748.3684
$ws.Range("J122").Value = 14499.6
$ws.Range("K122").Value = 8245.1052
$ws.Range("L122").Value = 43498.8
$ws.Range("M122").Value = -5795.1052
$ws.Range("N122").Value = -48398.8

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 7325.609
$ws.Range("I126").Value = 9559.532999999999
$ws.Range("J126").Value = 3137
$ws.Range("K126").Value = 28678.599
$ws.Range("L126").Value = 9411
$ws.Range("M126").Value = -26208.599
$ws.Range("N126").Value = -14351

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 900.1667
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
